# Updates TPM values for Fgf22-Fgfr2 LR-pair table and adds the new
# "Resolving-Mac" cluster (as both a sending and a target cluster).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
  @("FAPs", "Fgf22", "Fgfr2", "ECs", 1, 0.3333333333333333, 0.148113, 0.444339, 0.3768554877997065, 0.3768554877997065, 3, 1, 0.714474, 2.143422, 0.138796410342318, 0.138796410342318, 0.105822887562, 0.952405988058, 0.05230618892440249, 0.05230618892440249),
  @("FAPs", "Fgf22", "Fgfr2", "FAPs", 1, 0.3333333333333333, 0.148113, 0.444339, 0.3768554877997065, 0.3768554877997065, 3, 1, 4.140873, 12.422619, 0.8044215857867821, 0.8044215857867821, 0.613317122649, 5.519854103841, 0.3031506891082912, 0.3031506891082912),
  @("FAPs", "Fgf22", "Fgfr2", "MuSCs", 1, 0.3333333333333333, 0.148113, 0.444339, 0.3768554877997065, 0.3768554877997065, 3, 1, 0.2847646666666667, 0.8542940000000001, 0.05531945672713084, 0.05531945672713083, 0.042177349074, 0.379596141666, 0.02084744084971765, 0.02084744084971765),
  @("FAPs", "Fgf22", "Fgfr2", "Resolving-Mac", 1, 0.3333333333333333, 0.148113, 0.444339, 0.3768554877997065, 0.3768554877997065, 2, 0.6666666666666666, 0.007528666666666667, 0.022586, 0.00146254714376898, 0.00146254714376898, 0.001115093406, 0.010035840654, 0.0005511689172951265, 0.0005511689172951264),
  @("MuSCs", "Fgf22", "Fgfr2", "ECs", 1, 0.3333333333333333, 0.2449103333333333, 0.734731, 0.6231445122002934, 0.6231445122002934, 3, 1, 0.714474, 2.143422, 0.138796410342318, 0.138796410342318, 0.174982065498, 1.574838589482, 0.08649022141791553, 0.08649022141791553),
  @("MuSCs", "Fgf22", "Fgfr2", "FAPs", 1, 0.3333333333333333, 0.2449103333333333, 0.734731, 0.6231445122002934, 0.6231445122002934, 3, 1, 4.140873, 12.422619, 0.8044215857867821, 0.8044215857867821, 1.014142586721, 9.127283280489001, 0.5012708966784909, 0.5012708966784909),
  @("MuSCs", "Fgf22", "Fgfr2", "MuSCs", 1, 0.3333333333333333, 0.2449103333333333, 0.734731, 0.6231445122002934, 0.6231445122002934, 3, 1, 0.2847646666666667, 0.8542940000000001, 0.05531945672713084, 0.05531945672713083, 0.0697418094348889, 0.6276762849140001, 0.03447201587741319, 0.03447201587741319),
  @("MuSCs", "Fgf22", "Fgfr2", "Resolving-Mac", 1, 0.3333333333333333, 0.2449103333333333, 0.734731, 0.6231445122002934, 0.6231445122002934, 2, 0.6666666666666666, 0.007528666666666667, 0.022586, 0.00146254714376898, 0.00146254714376898, 0.001843848262888889, 0.016594634366, 0.0009113782264738535, 0.0009113782264738534)
)

$startRow = 2
for ($idx = 0; $idx -lt $rowsData.Count; $idx++) {
  $rowVals = $rowsData[$idx]
  $r = $startRow + $idx
  $arr = New-Object "object[,]" 1,20
  for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $rowVals[$i] }
  $rng = $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,20))
  $rng.Value = $arr
}

Write-Output "Updated rows 2-9"
